$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "generator_file" (sheet1): insert a new "prime_mover" column
# between generator_id and column_to_update, and append three new
# manual-correction rows (plant 56032 keep_leading_zeroes note, and
# the generator_id swap for plant 55350 generators 1/3).
# ------------------------------------------------------------------
$wsGen = $wb.Worksheets.Item("generator_file")

# Shift the existing "update" column (old D) out to the new E column,
# and the existing "column_to_update" column (old C) out to the new D
# column, freeing up column C for the new "prime_mover" data.
$wsGen.Range("D1:D5").Copy($wsGen.Range("E1:E5"))
$wsGen.Range("C1:C5").Copy($wsGen.Range("D1:D5"))
$wsGen.Range("C1:C5").ClearContents()

# Approximate column widths for the shifted/new columns (C, D, E).
$wsGen.Columns.Item(3).ColumnWidth = 17.6
$wsGen.Columns.Item(4).ColumnWidth = 16.8
$wsGen.Columns.Item(5).ColumnWidth = 14.8

# New header + values for column C ("prime_mover").
$wsGen.Range("C1").NumberFormat = "@"
$wsGen.Range("C1").Value = "prime_mover"
$wsGen.Range("C1").Font.Bold = $true

# New row 6: plant 56032 note about keeping leading zeroes.
$wsGen.Range("A6").Value = 56032
$wsGen.Range("D6").Value = "keep_leading_zeroes"
$wsGen.Range("D6").Style = "Normal"

# New row 7: plant 55350, generator 1 (CA) -> corrected generator_id 3.
$wsGen.Range("A7").Value = 55350
$wsGen.Range("B7").NumberFormat = "@"
$wsGen.Range("B7").Value = "1"
$wsGen.Range("C7").NumberFormat = "@"
$wsGen.Range("C7").Value = "CA"
$wsGen.Range("D7").Value = "generator_id"
$wsGen.Range("D7").Style = "Normal"
$wsGen.Range("E7").NumberFormat = "@"
$wsGen.Range("E7").Value = "3"

# New row 8: plant 55350, generator 3 (CT) -> corrected generator_id 1.
$wsGen.Range("A8").Value = 55350
$wsGen.Range("B8").NumberFormat = "@"
$wsGen.Range("B8").Value = "3"
$wsGen.Range("C8").NumberFormat = "@"
$wsGen.Range("C8").Value = "CT"
$wsGen.Range("D8").Value = "generator_id"
$wsGen.Range("D8").Style = "Normal"
$wsGen.Range("E8").NumberFormat = "@"
$wsGen.Range("E8").Value = "1"

# ------------------------------------------------------------------
# Update the saved cell selection on each sheet. "unit_file" is
# selected/activated last so it remains the active tab, matching the
# workbook's saved activeTab.
# ------------------------------------------------------------------
[void]$wsGen.Range("C2").Select()

$wsPlant = $wb.Worksheets.Item("plant_file")
[void]$wsPlant.Range("C14").Select()

$wsUnit = $wb.Worksheets.Item("unit_file")
[void]$wsUnit.Range("C8").Select()
